$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing Display-CF rows: amounts shown as "200.00" instead of "200" ---
$ws.Cells.Item(22,5).Value = "200.00"
$ws.Cells.Item(24,5).Value = "200.00"

# --- Column-letter -> number map used below ---
# --- New rows 25-33: additional "Display CF All Data" test cases ---
# Row 25
$ws.Cells.Item(25,1).Value = "Display CF All Data"
$ws.Cells.Item(25,3).Value = "24"
$ws.Cells.Item(25,4).Value = "2.7"
$ws.Cells.Item(25,5).Value = "5"
$ws.Cells.Item(25,6).Value = "26413421"
$ws.Cells.Item(25,7).Value = "PayNow"
$ws.Cells.Item(25,8).Value = "en_US"
$ws.Cells.Item(25,9).Value = "Elizath"
$ws.Cells.Item(25,10).Value = "Christine"
$ws.Cells.Item(25,11).Value = "258 Underwood rd"
$ws.Cells.Item(25,12).Value = "Suite 600"
$ws.Cells.Item(25,13).Value = "840"
$ws.Cells.Item(25,14).Value = "Arlington"
$ws.Cells.Item(25,15).Value = "VA"
$ws.Cells.Item(25,16).Value = "22201"
$ws.Cells.Item(25,18).Value = "Some Company"
$ws.Cells.Item(25,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(25,20).Value = "udf data 1"
$ws.Cells.Item(25,21).Value = "udf data 2"
$ws.Cells.Item(25,22).Value = "udf data 3"
$ws.Cells.Item(25,23).Value = "udf data 4"
$ws.Cells.Item(25,24).Value = "udf data 5"
$ws.Cells.Item(25,25).Value = "udf data 6"
$ws.Cells.Item(25,28).Value = "udf data 9"
$ws.Cells.Item(25,29).Value = "udf data 10"
$ws.Cells.Item(25,32).Value = "|Whole Wheat~$5| "
$ws.Rows.Item(25).RowHeight = 28.8

# Row 26
$ws.Cells.Item(26,1).Value = "Display CF All Data"
$ws.Cells.Item(26,3).Value = "25"
$ws.Cells.Item(26,4).Value = "2.7"
$ws.Cells.Item(26,5).Value = "200.00"
$ws.Cells.Item(26,6).Value = "26413421"
$ws.Cells.Item(26,7).Value = "PayNow"
$ws.Cells.Item(26,8).Value = "en_US"
$ws.Cells.Item(26,9).Value = "Elizath"
$ws.Cells.Item(26,10).Value = "Christine"
$ws.Cells.Item(26,11).Value = "258 Underwood rd"
$ws.Cells.Item(26,12).Value = "Suite 600"
$ws.Cells.Item(26,13).Value = "840"
$ws.Cells.Item(26,14).Value = "Arlington"
$ws.Cells.Item(26,15).Value = "VA"
$ws.Cells.Item(26,16).Value = "22201"
$ws.Cells.Item(26,18).Value = "Some Company"
$ws.Cells.Item(26,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(26,20).Value = "udf data 1"
$ws.Cells.Item(26,21).Value = "udf data 2"
$ws.Cells.Item(26,22).Value = "udf data 3"
$ws.Cells.Item(26,23).Value = "udf data 4"
$ws.Cells.Item(26,24).Value = "udf data 5"
$ws.Cells.Item(26,25).Value = "udf data 6"
$ws.Cells.Item(26,28).Value = "udf data 9"
$ws.Cells.Item(26,29).Value = "udf data 10"
$ws.Cells.Item(26,32).Value = "|Whole Wheat~$5| "
$ws.Rows.Item(26).RowHeight = 28.8

# Row 27
$ws.Cells.Item(27,1).Value = "Display CF All Data"
$ws.Cells.Item(27,3).Value = "26"
$ws.Cells.Item(27,4).Value = "3.0"
$ws.Cells.Item(27,5).Value = "5"
$ws.Cells.Item(27,6).Value = "26413422"
$ws.Cells.Item(27,7).Value = "PayNow"
$ws.Cells.Item(27,8).Value = "en_US"
$ws.Cells.Item(27,9).Value = "Elizath"
$ws.Cells.Item(27,10).Value = "Christine"
$ws.Cells.Item(27,11).Value = "258 Underwood rd"
$ws.Cells.Item(27,12).Value = "Suite 600"
$ws.Cells.Item(27,13).Value = "840"
$ws.Cells.Item(27,14).Value = "Arlington"
$ws.Cells.Item(27,15).Value = "VA"
$ws.Cells.Item(27,16).Value = "22201"
$ws.Cells.Item(27,18).Value = "Some Company"
$ws.Cells.Item(27,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(27,20).Value = "udf data 1"
$ws.Cells.Item(27,21).Value = "udf data 2"
$ws.Cells.Item(27,22).Value = "udf data 3"
$ws.Cells.Item(27,23).Value = "udf data 4"
$ws.Cells.Item(27,24).Value = "udf data 5"
$ws.Cells.Item(27,25).Value = "udf data 6"
$ws.Cells.Item(27,28).Value = "udf data 9"
$ws.Cells.Item(27,29).Value = "udf data 10"

# Row 28
$ws.Cells.Item(28,1).Value = "Display CF All Data"
$ws.Cells.Item(28,3).Value = "27"
$ws.Cells.Item(28,4).Value = "2.3"
$ws.Cells.Item(28,5).Value = "10.50"
$ws.Cells.Item(28,6).Value = "26413423"
$ws.Cells.Item(28,7).Value = "PayNow"
$ws.Cells.Item(28,8).Value = "en_US"
$ws.Cells.Item(28,9).Value = "Elizath"
$ws.Cells.Item(28,10).Value = "Christine"
$ws.Cells.Item(28,11).Value = "258 Underwood rd"
$ws.Cells.Item(28,12).Value = "Suite 600"
$ws.Cells.Item(28,13).Value = "840"
$ws.Cells.Item(28,14).Value = "Arlington"
$ws.Cells.Item(28,15).Value = "VA"
$ws.Cells.Item(28,16).Value = "22201"
$ws.Cells.Item(28,18).Value = "Some Company"
$ws.Cells.Item(28,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(28,20).Value = "udf data 1"
$ws.Cells.Item(28,21).Value = "udf data 2"
$ws.Cells.Item(28,22).Value = "udf data 3"
$ws.Cells.Item(28,23).Value = "udf data 4"
$ws.Cells.Item(28,24).Value = "udf data 5"
$ws.Cells.Item(28,25).Value = "udf data 6"
$ws.Cells.Item(28,28).Value = "udf data 9"
$ws.Cells.Item(28,29).Value = "udf data 10"

# Row 29
$ws.Cells.Item(29,1).Value = "Display CF All Data"
$ws.Cells.Item(29,3).Value = "28"
$ws.Cells.Item(29,4).Value = "2.5"
$ws.Cells.Item(29,5).Value = "10.50"
$ws.Cells.Item(29,6).Value = "26413424"
$ws.Cells.Item(29,7).Value = "PayNow"
$ws.Cells.Item(29,8).Value = "en_US"
$ws.Cells.Item(29,9).Value = "Elizath"
$ws.Cells.Item(29,10).Value = "Christine"
$ws.Cells.Item(29,11).Value = "258 Underwood rd"
$ws.Cells.Item(29,12).Value = "Suite 600"
$ws.Cells.Item(29,13).Value = "840"
$ws.Cells.Item(29,14).Value = "Arlington"
$ws.Cells.Item(29,15).Value = "VA"
$ws.Cells.Item(29,16).Value = "22201"
$ws.Cells.Item(29,18).Value = "Some Company"
$ws.Cells.Item(29,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(29,20).Value = "udf data 1"
$ws.Cells.Item(29,21).Value = "udf data 2"
$ws.Cells.Item(29,22).Value = "udf data 3"
$ws.Cells.Item(29,23).Value = "udf data 4"
$ws.Cells.Item(29,24).Value = "udf data 5"
$ws.Cells.Item(29,25).Value = "udf data 6"
$ws.Cells.Item(29,28).Value = "udf data 9"
$ws.Cells.Item(29,29).Value = "udf data 10"

# Row 30
$ws.Cells.Item(30,1).Value = "Display CF All Data"
$ws.Cells.Item(30,3).Value = "29"
$ws.Cells.Item(30,4).Value = "2.7"
$ws.Cells.Item(30,5).Value = "10.50"
$ws.Cells.Item(30,6).Value = "26413425"
$ws.Cells.Item(30,7).Value = "PayNow"
$ws.Cells.Item(30,8).Value = "en_US"
$ws.Cells.Item(30,9).Value = "Elizath"
$ws.Cells.Item(30,10).Value = "Christine"
$ws.Cells.Item(30,11).Value = "258 Underwood rd"
$ws.Cells.Item(30,12).Value = "Suite 600"
$ws.Cells.Item(30,13).Value = "840"
$ws.Cells.Item(30,14).Value = "Arlington"
$ws.Cells.Item(30,15).Value = "VA"
$ws.Cells.Item(30,16).Value = "22201"
$ws.Cells.Item(30,18).Value = "Some Company"
$ws.Cells.Item(30,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(30,20).Value = "udf data 1"
$ws.Cells.Item(30,21).Value = "udf data 2"
$ws.Cells.Item(30,22).Value = "udf data 3"
$ws.Cells.Item(30,23).Value = "udf data 4"
$ws.Cells.Item(30,24).Value = "udf data 5"
$ws.Cells.Item(30,25).Value = "udf data 6"
$ws.Cells.Item(30,28).Value = "udf data 9"
$ws.Cells.Item(30,29).Value = "udf data 10"
$ws.Cells.Item(30,32).Value = "|Whole Wheat~$5| "
$ws.Rows.Item(30).RowHeight = 28.8

# Row 31
$ws.Cells.Item(31,1).Value = "Display CF All Data"
$ws.Cells.Item(31,3).Value = "30"
$ws.Cells.Item(31,4).Value = "3.0"
$ws.Cells.Item(31,5).Value = "10.50"
$ws.Cells.Item(31,6).Value = "26413426"
$ws.Cells.Item(31,7).Value = "PayNow"
$ws.Cells.Item(31,8).Value = "en_US"
$ws.Cells.Item(31,9).Value = "Elizath"
$ws.Cells.Item(31,10).Value = "Christine"
$ws.Cells.Item(31,11).Value = "258 Underwood rd"
$ws.Cells.Item(31,12).Value = "Suite 600"
$ws.Cells.Item(31,13).Value = "840"
$ws.Cells.Item(31,14).Value = "Arlington"
$ws.Cells.Item(31,15).Value = "VA"
$ws.Cells.Item(31,16).Value = "22201"
$ws.Cells.Item(31,18).Value = "Some Company"
$ws.Cells.Item(31,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(31,20).Value = "udf data 1"
$ws.Cells.Item(31,21).Value = "udf data 2"
$ws.Cells.Item(31,22).Value = "udf data 3"
$ws.Cells.Item(31,23).Value = "udf data 4"
$ws.Cells.Item(31,24).Value = "udf data 5"
$ws.Cells.Item(31,25).Value = "udf data 6"
$ws.Cells.Item(31,28).Value = "udf data 9"
$ws.Cells.Item(31,29).Value = "udf data 10"

# Row 32
$ws.Cells.Item(32,1).Value = "Display CF All Data"
$ws.Cells.Item(32,3).Value = "31"
$ws.Cells.Item(32,4).Value = "2.3"
$ws.Cells.Item(32,5).Value = "10.899"
$ws.Cells.Item(32,6).Value = "26413427"
$ws.Cells.Item(32,7).Value = "PayNow"
$ws.Cells.Item(32,8).Value = "en_US"
$ws.Cells.Item(32,9).Value = "Elizath"
$ws.Cells.Item(32,10).Value = "Christine"
$ws.Cells.Item(32,11).Value = "258 Underwood rd"
$ws.Cells.Item(32,12).Value = "Suite 600"
$ws.Cells.Item(32,13).Value = "840"
$ws.Cells.Item(32,14).Value = "Arlington"
$ws.Cells.Item(32,15).Value = "VA"
$ws.Cells.Item(32,16).Value = "22201"
$ws.Cells.Item(32,18).Value = "Some Company"
$ws.Cells.Item(32,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(32,20).Value = "udf data 1"
$ws.Cells.Item(32,21).Value = "udf data 2"
$ws.Cells.Item(32,22).Value = "udf data 3"
$ws.Cells.Item(32,23).Value = "udf data 4"
$ws.Cells.Item(32,24).Value = "udf data 5"
$ws.Cells.Item(32,25).Value = "udf data 6"
$ws.Cells.Item(32,28).Value = "udf data 9"
$ws.Cells.Item(32,29).Value = "udf data 10"

# Row 33
$ws.Cells.Item(33,1).Value = "Display CF All Data"
$ws.Cells.Item(33,3).Value = "32"
$ws.Cells.Item(33,4).Value = "2.3"
$ws.Cells.Item(33,5).Value = "20.8899"
$ws.Cells.Item(33,6).Value = "26413428"
$ws.Cells.Item(33,7).Value = "PayNow"
$ws.Cells.Item(33,8).Value = "en_US"
$ws.Cells.Item(33,9).Value = "Elizath"
$ws.Cells.Item(33,10).Value = "Christine"
$ws.Cells.Item(33,11).Value = "258 Underwood rd"
$ws.Cells.Item(33,12).Value = "Suite 600"
$ws.Cells.Item(33,13).Value = "840"
$ws.Cells.Item(33,14).Value = "Arlington"
$ws.Cells.Item(33,15).Value = "VA"
$ws.Cells.Item(33,16).Value = "22201"
$ws.Cells.Item(33,18).Value = "Some Company"
$ws.Cells.Item(33,19).Value = "iahmed@govolution.com"
$ws.Cells.Item(33,20).Value = "udf data 1"
$ws.Cells.Item(33,21).Value = "udf data 2"
$ws.Cells.Item(33,22).Value = "udf data 3"
$ws.Cells.Item(33,23).Value = "udf data 4"
$ws.Cells.Item(33,24).Value = "udf data 5"
$ws.Cells.Item(33,25).Value = "udf data 6"
$ws.Cells.Item(33,28).Value = "udf data 9"
$ws.Cells.Item(33,29).Value = "udf data 10"

# --- Column S needs the bordered "email" style (style index 2) like rows 2-24 ---
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(25,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(26,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(27,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(28,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(29,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(30,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(31,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(32,19).PasteSpecial(-4122)
$ws.Cells.Item(21,19).Copy()
$ws.Cells.Item(33,19).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view: scroll down and select F33 (matches the authored sheetView) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("F33").Select()
